# ICS496-Fall24-FinalPoster.pptx - grammar/typo fixes
#   - Add missing period to the "Challenges" bullet about defining issues/scope.
#   - Add missing period to the "Challenges" bullet about documentation.
#   - Fix "mailist list" -> "mailing list" typo in the Accomplishments section.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 144 ("Challenges") holds both of the first two bullet paragraphs,
# each as a single run - update the run text directly so PowerPoint keeps
# a single <a:r> per paragraph instead of splitting off the appended text.
$challengesShape = $s.Shapes.Item(11)
$challengesText = $challengesShape.TextFrame.TextRange

$challengesText.Paragraphs(1, 1).Runs(1, 1).Text = "Defining issues and scope according to specifications, and refining solutions through iterative feedback."
$challengesText.Paragraphs(3, 1).Runs(1, 1).Text = "Developing clear and concise documentation, providing guidance to users potentially unfamiliar with the technology."

# Shape 152 ("Accomplishments") - the "mailist list" typo is the second run
# of its 8th paragraph ("Implemented download and refresh of " + the typo run).
$accomplishmentsShape = $s.Shapes.Item(19)
$accomplishmentsText = $accomplishmentsShape.TextFrame.TextRange

$accomplishmentsText.Paragraphs(8, 1).Runs(2, 1).Text = "two mailing list data archives (mod mbox and pipermail)."
